$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above current row 6, shifting existing rows 6-46 down to 7-47.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new reference entry.
$ws.Cells.Item(6, 1).Value = "Deter et al"
$ws.Cells.Item(6, 2).Value = "Antibiotic tolerance is associated with a broad and complex transcriptional response in E. coli"
$ws.Cells.Item(6, 3).Value = "10.1038/s41598-021-85509-7"
$ws.Cells.Item(6, 4).Value = "Ampicillin resistance"
$ws.Cells.Item(6, 5).Value = "Deter et al. generated RNA-seq data on both antibiotic-treated and -untreated populations emerging from stationary phase."

# Reflect the cursor/selection position recorded in the saved workbook.
$ws.Range("C30").Select()
